$d = $word.ActiveDocument

# Highlight color used for quantitative impact metrics (matches target hex 2C3E50).
# Word/COM Font.Color uses BGR-packed integers (0x00BBGGRR), not RGB.
function Get-WordColorFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}
$metricColor = Get-WordColorFromHex "2C3E50"

# Applies bold + the metric color to every occurrence of $token found inside the
# paragraph at index $paraIndex (1-based), scanning left-to-right so repeated
# tokens within the same paragraph (e.g. two different percentages) are each
# located independently without drifting into neighboring paragraphs.
function Set-MetricHighlight($paraIndex, $token) {
    $para = $d.Paragraphs.Item($paraIndex)
    $paraStart = $para.Range.Start
    $paraEnd = $para.Range.End

    if ($para.Range.Text -notlike "*$token*") {
        Write-Output "WARNING: paragraph $paraIndex does not contain '$token' - skipping"
        return
    }

    $searchStart = $paraStart
    while ($true) {
        $rng = $d.Range($searchStart, $paraEnd)
        $found = $rng.Find.Execute($token, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        if ($rng.Start -lt $paraStart -or $rng.End -gt $paraEnd) { break }

        $rng.Font.Bold = $true
        $rng.Font.Color = $metricColor

        $searchStart = $rng.End
        if ($searchStart -ge $paraEnd) { break }
    }
}

# Paragraph 9: "• Discovered systematic race coding errors ... from 23% to 64%"
Set-MetricHighlight 9 "23%"
Set-MetricHighlight 9 "64%"

# Paragraph 11: "• Achieved 87% ... of 71%, reducing polling error margins from ±4.2% to ±2.1%"
Set-MetricHighlight 11 "87%"
Set-MetricHighlight 11 "71%"
Set-MetricHighlight 11 "±4.2%"
Set-MetricHighlight 11 "±2.1%"

# Paragraph 31: "• Wrote RFP and analyzed bids from 1,200 vendors ..."
Set-MetricHighlight 31 "1,200"

# Paragraph 46: "• Created comprehensive meta-analysis ... $400M ... now valued at $1B+"
Set-MetricHighlight 46 "$400M"
Set-MetricHighlight 46 "$1B"

# Paragraph 63: "• Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Set-MetricHighlight 63 "73.5%"
Set-MetricHighlight 63 "$4.7M"

# Paragraph 65: "• Achieved 87% prediction accuracy ... of 71%" (short variant)
Set-MetricHighlight 65 "87%"
Set-MetricHighlight 65 "71%"

Write-Output "Done applying metric highlighting"
